# Bump the "Förändrad" (Changed) date in column C for every data row
# (rows 2-536) from Excel serial date 45203 (2023-10-04) to
# 45204 (2023-10-05), i.e. +1 day, leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column C (falls back to 536 which is
# the known extent of the data in this workbook).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 536 }

$rng = $ws.Range("C2:C$lastRow")
$rng.Value2 = 45204
